# Fixed update to excel issue
#
# 1) Rename the "Requested quantity" header on the existing sheets.
# 2) Add a new "PO Forecast" sheet (after the existing sheets) with a
#    ds / PO_Forecast / yhat_lower / yhat_upper table.

$wb = $excel.ActiveWorkbook

# --- 1) Rename headers on existing sheets -------------------------------
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- 2) Add the new "PO Forecast" sheet at the end ----------------------
$wsForecast = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsForecast.Name = "PO Forecast"

# Header row
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Match the bold/centered/bordered header style used on the other sheets
$wsWeekly.Range("A1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)

$data = @(
    @(44990.99999999999, 5, -30.14174783691047, 43.49008578149358),
    @(44997.99999999999, 6, -30.15114135647556, 43.23571565869473),
    @(45025.99999999999, 8, -29.64767195016909, 46.77194602210746),
    @(45032.99999999999, 8, -29.27039547496476, 44.71385591691231),
    @(45039.99999999999, 9, -25.29331266598268, 45.26168244002804),
    @(45046.99999999999, 10, -25.6137617201564, 45.6721514526367),
    @(45053.99999999999, 10, -26.37891129306032, 47.86766718887141),
    @(45060.99999999999, 11, -25.1802436057897, 47.51268357002152),
    @(45074.99999999999, 12, -26.07879685596586, 44.91237685261807),
    @(45088.99999999999, 13, -24.31688992415014, 49.03298799428363),
    @(45095.99999999999, 14, -23.51823798933276, 49.88217786435019),
    @(45116.99999999999, 15, -20.62004365820463, 49.84753721973271),
    @(45123.99999999999, 16, -20.46582299364672, 50.94495173332315),
    @(45130.99999999999, 16, -20.00922264123448, 53.71330279882361),
    @(45144.99999999999, 18, -16.52121239304573, 53.35880680795664),
    @(45151.99999999999, 18, -17.42096575427762, 53.33142158820787),
    @(45228.99999999999, 24, -9.947147556420024, 61.19927244246673),
    @(45389.99999999999, 37, 1.785274799117118, 73.53673762917946),
    @(45396.99999999999, 38, 2.749665781297343, 71.95250398005794),
    @(45410.99999999999, 39, 3.269394181356555, 73.61735354564229),
    @(45417.99999999999, 40, 6.216926347565171, 74.21100818927178),
    @(45424.99999999999, 40, 4.17550970416975, 76.55500892672045),
    @(45431.99999999999, 41, 5.449114476607325, 76.15094560552511),
    @(45438.99999999999, 41, 3.478222087653547, 73.64096076918105),
    @(45445.99999999999, 42, 7.024958098288044, 77.87006727251836),
    @(45480.99999999999, 45, 9.370017440537216, 80.40524873535826),
    @(45487.99999999999, 45, 8.438988761532592, 80.95187370487066),
    @(45494.99999999999, 46, 9.211500289929095, 79.48014067911861),
    @(45508.99999999999, 47, 10.32806007639434, 81.79874419784041),
    @(45515.99999999999, 48, 12.51898749926939, 84.57657478943156),
    @(45543.99999999999, 50, 15.95572653619961, 85.79351076418571),
    @(45550.99999999999, 50, 14.22928568448568, 85.05811037087186),
    @(45557.99999999999, 51, 17.72213979009009, 85.31792648631719),
    @(45564.99999999999, 51, 15.07657425636121, 89.03686803134033),
    @(45571.99999999999, 52, 17.9698755586979, 89.10540379455993),
    @(45578.99999999999, 53, 15.31284531765765, 87.41987604099133),
    @(45585.99999999999, 53, 19.35377040612568, 86.88208352573973),
    @(45592.99999999999, 54, 16.31589188584618, 87.16925932189604),
    @(45599.99999999999, 54, 16.65714945186961, 89.71322982464913)
)

$r = 2
foreach ($row in $data) {
    $wsForecast.Cells.Item($r, 1).Value = $row[0]
    $wsForecast.Cells.Item($r, 2).Value = $row[1]
    $wsForecast.Cells.Item($r, 3).Value = $row[2]
    $wsForecast.Cells.Item($r, 4).Value = $row[3]
    $r++
}

# Match the date-formatted style used in column A of the other sheets
$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A40").PasteSpecial(-4122)

# Re-apply the values after the format paste (PasteSpecial formats only,
# but keep this for safety / idempotency)
$r = 2
foreach ($row in $data) {
    $wsForecast.Cells.Item($r, 1).Value = $row[0]
    $r++
}

# Keep the originally-active sheet/selection (only the sheet list changed)
$wsWeekly.Activate() | Out-Null
$wsWeekly.Range("A1").Select() | Out-Null

Write-Output "PO Forecast sheet added; headers renamed."
